$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so Excel does not
# reinterpret them as numbers (these cells are inline strings in the source).
# NumberFormat/Style must be applied per-cell (not as a multi-area union).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.139.99"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "1.805.29"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "338.60"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.3935"
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("D8").Value = "0.3492"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "0.07535"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "22.03"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "6.511"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "1.805.58"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "7.149"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "0.00001102"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "0.06717"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "85.06"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").Value = "6.569"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "28.143.02"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").Value = "12.44"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "2.412"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "1.490"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "21.41"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "2.519"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "154.08"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "2.012.56"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "135.48"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "6.182"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "4.017"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").Value = "0.08843"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").Value = "13.10"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "0.6944"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "0.06547"
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02422"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.442"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.606"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "0.2215"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "1.255"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "8.483"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").Value = "14.60"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "0.6421"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "2.147"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "131.08"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "0.07184"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "80.12"
$ws.Range("E51").Value = "  +0.91%  "

# Restore default cell style so the saved XML has no stray "s" attribute
# on these cells (matches the original, unstyled inline-string cells).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

